# Commit: "changed binary cost files to €/g"
# The raw material cost in economical_params!B2 was stored in €/kg (34.9)
# and is converted to €/g (0.035). The Scaling sheet formulas
# (B2 = economical_params!B2*0.8, C2 = economical_params!B2*2) are left
# untouched and simply recompute from the new base value.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("economical_params")

# Convert the virgin-PLA cost from €/kg to €/g.
$ws1.Range("B2").Value = 0.035

# The author's last interaction was on economical_params (cell B17),
# making it the active sheet/tab instead of Scaling.
$ws1.Activate()
$ws1.Range("B17").Select()
